$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row (row 3) - mirrors the structure of row 2, new study entry
$ws.Range("A3").Value = 4564654
$ws.Range("B3").Value = 64654
$ws.Range("C3").Value = "BRANY"
$ws.Range("D3").Value = 123465498
$ws.Range("E3").Value = "Einstein"
$ws.Range("F3").Value = "A study"
$ws.Range("G3").Value = "A study about bla bla bla"
$ws.Range("H3").Value = $true
$ws.Range("K3").Value = "ad"
$ws.Range("Q3").Value = $false
$ws.Range("S3").Value = "McDonald"
$ws.Range("T3").Value = "Mike"
$ws.Range("U3").Value = "John"
$ws.Range("W3").Value = 64554

# Hyperlink cells (email + document link), matching row 2's pattern
$ws.Hyperlinks.Add($ws.Range("V3"), "mailto:asd@sdf.dd", "", "", "asd@sdf.dd")
$ws.Range("V3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("X3"), "http://www.google.com/", "", "", "www.google.com")
$ws.Range("X3").Style = "Hyperlink"

# Column D widened to fit the longer IRBNumber value
$ws.Columns.Item(4).ColumnWidth = 9.14

# Leave selection on the newly entered row, as Excel would after data entry
$ws.Range("A3").Select() | Out-Null
